$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text columns (G:K) from auto-number conversion
$ws.Range("G3:K12").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = ' Oct 18 2020'
$ws.Range("B3").Value = ' Dubai (DSC)'
$ws.Range("C3").Value = 'Match tied (Kings XI won the one-over eliminator)'
$ws.Range("D3").Value = 'Mumbai Indians'
$ws.Range("E3").Value = 'Kings XI Punjab'
$ws.Range("F3").Value = 'Kieron Pollard '
$ws.Range("G3").Value = '34'
$ws.Range("H3").Value = '12'
$ws.Range("I3").Value = '1'
$ws.Range("J3").Value = '4'
$ws.Range("K3").Value = '283.33'

# Row 4
$ws.Range("A4").Value = ' Oct 25 2020'
$ws.Range("B4").Value = ' Abu Dhabi'
$ws.Range("C4").Value = 'Royals won by 8 wickets (with 10 balls remaining)'
$ws.Range("D4").Value = 'Mumbai Indians'
$ws.Range("E4").Value = 'Rajasthan Royals'
$ws.Range("F4").Value = 'Kieron Pollard '
$ws.Range("G4").Value = '6'
$ws.Range("H4").Value = '4'
$ws.Range("I4").Value = '0'
$ws.Range("J4").Value = '1'
$ws.Range("K4").Value = '150.00'

# Row 5
$ws.Range("A5").Value = ' Oct 4 2020'
$ws.Range("B5").Value = ' Sharjah'
$ws.Range("C5").Value = 'Mumbai won by 34 runs'
$ws.Range("D5").Value = 'Mumbai Indians'
$ws.Range("E5").Value = 'Sunrisers Hyderabad'
$ws.Range("F5").Value = 'Kieron Pollard '
$ws.Range("G5").Value = '25'
$ws.Range("H5").Value = '13'
$ws.Range("I5").Value = '0'
$ws.Range("J5").Value = '3'
$ws.Range("K5").Value = '192.30'

# Row 6
$ws.Range("A6").Value = ' Oct 11 2020'
$ws.Range("B6").Value = ' Abu Dhabi'
$ws.Range("C6").Value = 'Mumbai won by 5 wickets (with 2 balls remaining)'
$ws.Range("D6").Value = 'Mumbai Indians'
$ws.Range("E6").Value = 'Delhi Capitals'
$ws.Range("F6").Value = 'Kieron Pollard '
$ws.Range("G6").Value = '11'
$ws.Range("H6").Value = '14'
$ws.Range("I6").Value = '1'
$ws.Range("J6").Value = '0'
$ws.Range("K6").Value = '78.57'

# Row 7
$ws.Range("A7").Value = ' Nov 5 2020'
$ws.Range("B7").Value = ' Dubai (DSC)'
$ws.Range("C7").Value = 'Mumbai won by 57 runs'
$ws.Range("D7").Value = 'Mumbai Indians'
$ws.Range("E7").Value = 'Delhi Capitals'
$ws.Range("F7").Value = 'Kieron Pollard '
$ws.Range("G7").Value = '0'
$ws.Range("H7").Value = '2'
$ws.Range("I7").Value = '0'
$ws.Range("J7").Value = '0'
$ws.Range("K7").Value = '0.00'

# Row 8
$ws.Range("A8").Value = ' Sep 19 2020'
$ws.Range("B8").Value = ' Abu Dhabi'
$ws.Range("C8").Value = 'Super Kings won by 5 wickets (with 4 balls remaining)'
$ws.Range("D8").Value = 'Mumbai Indians'
$ws.Range("E8").Value = 'Chennai Super Kings'
$ws.Range("F8").Value = 'Kieron Pollard '
$ws.Range("G8").Value = '18'
$ws.Range("H8").Value = '14'
$ws.Range("I8").Value = '1'
$ws.Range("J8").Value = '1'
$ws.Range("K8").Value = '128.57'

# Row 9
$ws.Range("A9").Value = ' Nov 3 2020'
$ws.Range("B9").Value = ' Sharjah'
$ws.Range("C9").Value = 'Sunrisers won by 10 wickets (with 17 balls remaining)'
$ws.Range("D9").Value = 'Mumbai Indians'
$ws.Range("E9").Value = 'Sunrisers Hyderabad'
$ws.Range("F9").Value = 'Kieron Pollard '
$ws.Range("G9").Value = '41'
$ws.Range("H9").Value = '25'
$ws.Range("I9").Value = '2'
$ws.Range("J9").Value = '4'
$ws.Range("K9").Value = '164.00'

# Row 10
$ws.Range("A10").Value = ' Sep 28 2020'
$ws.Range("B10").Value = ' Dubai (DSC)'
$ws.Range("C10").Value = 'Match tied (RCB won the one-over eliminator)'
$ws.Range("D10").Value = 'Mumbai Indians'
$ws.Range("E10").Value = 'Royal Challengers Bangalore'
$ws.Range("F10").Value = 'Kieron Pollard '
$ws.Range("G10").Value = '60'
$ws.Range("H10").Value = '24'
$ws.Range("I10").Value = '3'
$ws.Range("J10").Value = '5'
$ws.Range("K10").Value = '250.00'

# Row 11
$ws.Range("A11").Value = ' Oct 1 2020'
$ws.Range("B11").Value = ' Abu Dhabi'
$ws.Range("C11").Value = 'Mumbai won by 48 runs'
$ws.Range("D11").Value = 'Mumbai Indians'
$ws.Range("E11").Value = 'Kings XI Punjab'
$ws.Range("F11").Value = 'Kieron Pollard '
$ws.Range("G11").Value = '47'
$ws.Range("H11").Value = '20'
$ws.Range("I11").Value = '3'
$ws.Range("J11").Value = '4'
$ws.Range("K11").Value = '235.00'

# Row 12
$ws.Range("A12").Value = ' Sep 23 2020'
$ws.Range("B12").Value = ' Abu Dhabi'
$ws.Range("C12").Value = 'Mumbai won by 49 runs'
$ws.Range("D12").Value = 'Mumbai Indians'
$ws.Range("E12").Value = 'Kolkata Knight Riders'
$ws.Range("F12").Value = 'Kieron Pollard '
$ws.Range("G12").Value = '13'
$ws.Range("H12").Value = '7'
$ws.Range("I12").Value = '1'
$ws.Range("J12").Value = '0'
$ws.Range("K12").Value = '185.71'
